# Auto-generated Excel COM-interop script
# Applies scheduled market-data refresh updates to the Leve profit sheets
# (columns H:N = currentAveragePrice/NQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)

$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("H62").Value = 1925
$ws.Range("I62").Value = 1925
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1925
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = $null
$ws.Range("N62").Value = -1301
$ws.Range("H65").Value = 1925
$ws.Range("I65").Value = 1925
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 9625
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = $null
$ws.Range("N65").Value = -6505
$ws.Range("H70").Value = 1229.5883
$ws.Range("I70").Value = 1666.8334
$ws.Range("J70").Value = 991.0909
$ws.Range("K70").Value = 5000.5002
$ws.Range("L70").Value = 2973.2727
$ws.Range("M70").Value = -4730.5002
$ws.Range("N70").Value = -3513.2727
$ws.Range("H73").Value = 1229.5883
$ws.Range("I73").Value = 1666.8334
$ws.Range("J73").Value = 991.0909
$ws.Range("K73").Value = 5000.5002
$ws.Range("L73").Value = 2973.2727
$ws.Range("M73").Value = -4064.5002
$ws.Range("N73").Value = -4845.2727
$ws.Range("H76").Value = 3409
$ws.Range("J76").Value = 4800
$ws.Range("L76").Value = 4800
$ws.Range("N76").Value = -5430
$ws.Range("H79").Value = 3409
$ws.Range("J79").Value = 4800
$ws.Range("L79").Value = 4800
$ws.Range("N79").Value = -6984
$ws.Range("H98").Value = 1068
$ws.Range("I98").Value = 853.3333
$ws.Range("J98").Value = 3000
$ws.Range("K98").Value = 853.3333
$ws.Range("L98").Value = 3000
$ws.Range("M98").Value = 644.6667
$ws.Range("N98").Value = -5996
$ws.Range("H122").Value = 1068
$ws.Range("I122").Value = 853.3333
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 2559.9999
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -109.9998999999998
$ws.Range("N122").Value = -13900
$ws.Range("H129").Value = 834.12164
$ws.Range("I129").Value = 462.1111
$ws.Range("J129").Value = 885.6308
$ws.Range("K129").Value = 1386.3333
$ws.Range("L129").Value = 2656.8924
$ws.Range("M129").Value = 3613.6667
$ws.Range("N129").Value = -12656.8924
$ws.Range("H137").Value = 1097.9259
$ws.Range("I137").Value = 968.7
$ws.Range("J137").Value = 1467.1428
$ws.Range("K137").Value = 2906.1
$ws.Range("L137").Value = 4401.428400000001
$ws.Range("M137").Value = -356.1000000000004
$ws.Range("N137").Value = -9501.428400000001
# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("H32").Value = 38294.21
$ws.Range("J32").Value = 155821.2
$ws.Range("L32").Value = 155821.2
$ws.Range("N32").Value = -156395.2
$ws.Range("H63").Value = 2618.1765
$ws.Range("J63").Value = 3316.6667
$ws.Range("L63").Value = 3316.6667
$ws.Range("N63").Value = -4688.6667
$ws.Range("H66").Value = 2618.1765
$ws.Range("J66").Value = 3316.6667
$ws.Range("L66").Value = 16583.3335
$ws.Range("N66").Value = -23447.3335
# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("H62").Value = 47984.75
$ws.Range("J62").Value = 47979.668
$ws.Range("L62").Value = 47979.668
$ws.Range("N62").Value = -49351.668
$ws.Range("H65").Value = 47984.75
$ws.Range("J65").Value = 47979.668
$ws.Range("L65").Value = 143939.004
$ws.Range("N65").Value = -150803.004
$ws.Range("H99").Value = 1563.3334
$ws.Range("I99").Value = 1473.3334
$ws.Range("J99").Value = 1833.3334
$ws.Range("K99").Value = 1473.3334
$ws.Range("L99").Value = 1833.3334
$ws.Range("M99").Value = 24.66660000000002
$ws.Range("N99").Value = -4829.3334
$ws.Range("H105").Value = 201808.9
$ws.Range("I105").Value = 144897
$ws.Range("J105").Value = 334603.34
$ws.Range("K105").Value = 144897
$ws.Range("L105").Value = 334603.34
$ws.Range("M105").Value = -143150
$ws.Range("N105").Value = -338097.34
$ws.Range("H107").Value = 37039960
$ws.Range("I107").Value = 66669130
$ws.Range("K107").Value = 66669130
$ws.Range("M107").Value = -66667210
# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("H58").Value = 7678.325
$ws.Range("I58").Value = 1263.6666
$ws.Range("K58").Value = 1263.6666
$ws.Range("M58").Value = -1060.6666
$ws.Range("H62").Value = 4632061.5
$ws.Range("I62").Value = 55555556
$ws.Range("J62").Value = 2652.818
$ws.Range("K62").Value = 55555556
$ws.Range("L62").Value = 2652.818
$ws.Range("M62").Value = -55554932
$ws.Range("N62").Value = -3900.818
$ws.Range("H63").Value = 40720
$ws.Range("J63").Value = 40720
$ws.Range("L63").Value = 40720
$ws.Range("N63").Value = -42092
$ws.Range("H65").Value = 4632061.5
$ws.Range("I65").Value = 55555556
$ws.Range("J65").Value = 2652.818
$ws.Range("K65").Value = 277777780
$ws.Range("L65").Value = 13264.09
$ws.Range("M65").Value = -277774660
$ws.Range("N65").Value = -19504.09
$ws.Range("H66").Value = 40720
$ws.Range("J66").Value = 40720
$ws.Range("L66").Value = 122160
$ws.Range("N66").Value = -129024
$ws.Range("H80").Value = 10337
$ws.Range("I80").Value = 5000
$ws.Range("J80").Value = 11099.429
$ws.Range("K80").Value = 5000
$ws.Range("L80").Value = 11099.429
$ws.Range("M80").Value = -3877
$ws.Range("N80").Value = -13345.429
$ws.Range("H83").Value = 10337
$ws.Range("I83").Value = 5000
$ws.Range("J83").Value = 11099.429
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 33298.287
$ws.Range("M83").Value = -9384
$ws.Range("N83").Value = -44530.287
$ws.Range("H136").Value = 7678.325
$ws.Range("I136").Value = 1263.6666
$ws.Range("K136").Value = 3790.9998
$ws.Range("M136").Value = -1240.9998
# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("H98").Value = 53384.285
$ws.Range("I98").Value = 487.875
$ws.Range("J98").Value = 85935.92
$ws.Range("K98").Value = 1463.625
$ws.Range("L98").Value = 257807.76
$ws.Range("M98").Value = 34.375
$ws.Range("N98").Value = -260803.76
$ws.Range("H131").Value = 791.1799999999999
$ws.Range("J131").Value = 819
$ws.Range("L131").Value = 2457
$ws.Range("N131").Value = -12537
# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("H64").Value = 46551.5
$ws.Range("J64").Value = 46551.5
$ws.Range("L64").Value = 46551.5
$ws.Range("N64").Value = -47047.5
$ws.Range("H67").Value = 46551.5
$ws.Range("J67").Value = 46551.5
$ws.Range("L67").Value = 46551.5
$ws.Range("N67").Value = -48267.5
$ws.Range("H102").Value = 752999
$ws.Range("I102").Value = 3427.4285
$ws.Range("K102").Value = 3427.4285
$ws.Range("M102").Value = -1805.4285
$ws.Range("H113").Value = 1854.8889
$ws.Range("J113").Value = 1999.1428
$ws.Range("L113").Value = 1999.1428
$ws.Range("N113").Value = -6339.1428
$ws.Range("H141").Value = 37100
$ws.Range("J141").Value = 37100
$ws.Range("L141").Value = 37100
$ws.Range("N141").Value = -47460
# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("H62").Value = 43248
$ws.Range("J62").Value = 43248
$ws.Range("L62").Value = 43248
$ws.Range("N62").Value = -44496
$ws.Range("H65").Value = 43248
$ws.Range("J65").Value = 43248
$ws.Range("L65").Value = 129744
$ws.Range("N65").Value = -135984
$ws.Range("H68").Value = 3610.6924
$ws.Range("I68").Value = 1500.5
$ws.Range("J68").Value = 4548.5557
$ws.Range("K68").Value = 1500.5
$ws.Range("L68").Value = 4548.5557
$ws.Range("M68").Value = -751.5
$ws.Range("N68").Value = -6046.5557
$ws.Range("H71").Value = 3610.6924
$ws.Range("I71").Value = 1500.5
$ws.Range("J71").Value = 4548.5557
$ws.Range("K71").Value = 7502.5
$ws.Range("L71").Value = 22742.7785
$ws.Range("M71").Value = -3758.5
$ws.Range("N71").Value = -30230.7785
$ws.Range("H122").Value = 2004.5
$ws.Range("J122").Value = 2004.5
$ws.Range("L122").Value = 6013.5
$ws.Range("N122").Value = -10913.5
$ws.Range("H132").Value = 4090.2173
$ws.Range("I132").Value = 4829.375
$ws.Range("J132").Value = 2400.7144
$ws.Range("K132").Value = 14488.125
$ws.Range("L132").Value = 7202.1432
$ws.Range("M132").Value = -11958.125
$ws.Range("N132").Value = -12262.1432
# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("H122").Value = 1550.4375
$ws.Range("I122").Value = 1258.6666
$ws.Range("J122").Value = 2425.75
$ws.Range("K122").Value = 3775.9998
$ws.Range("L122").Value = 7277.25
$ws.Range("M122").Value = -1325.9998
$ws.Range("N122").Value = -12177.25
